# Generate Report for Handoff
# Updates status cells from "In Translation" to "Ready for handoff" and
# refreshes the associated timestamp cells, on all three report sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: per-language status columns (E, F) + generate-date column (G)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-16 20:35:56"

# zh-cn sheet: Status column (C) + Latest Handoff Datetime column (H)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-16 20:35:51"

# de-de sheet: Status column (C) only; its Latest Handoff Datetime (H2)
# shares the same underlying text as Overview!G2 and is refreshed to match.
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-16 20:35:56"
